# Corrects mixed-up BT/BLE power-state labels and refreshes the
# measured current-draw stats that were captured under the wrong labels
# because of a data-collection threading race; also drops the stray
# extra row that the race condition produced.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
# --- Sheet "3_3": relabel rows and fix values ---
$ws1.Cells.Item(2, 1).Value = "Deep_Sleep"
$ws1.Cells.Item(2, 2).Value = 0.0994
$ws1.Cells.Item(2, 3).Value = 0.1334
$ws1.Cells.Item(2, 4).Value = 0.075
$ws1.Cells.Item(2, 5).Value = 0.0179
$ws1.Cells.Item(2, 6).Value = 100
$ws1.Cells.Item(3, 1).Value = "BT_Idle"
$ws1.Cells.Item(3, 2).Value = 5.4892
$ws1.Cells.Item(3, 3).Value = 5.4989
$ws1.Cells.Item(3, 4).Value = 5.4768
$ws1.Cells.Item(3, 5).Value = 0.0045
$ws1.Cells.Item(3, 6).Value = 100
$ws1.Cells.Item(4, 1).Value = "BT_Pscan"
$ws1.Cells.Item(4, 2).Value = 0.1454
$ws1.Cells.Item(4, 3).Value = 1.3266
$ws1.Cells.Item(4, 4).Value = 0.0623
$ws1.Cells.Item(4, 5).Value = 0.2012
$ws1.Cells.Item(4, 6).Value = 100
$ws1.Cells.Item(5, 1).Value = "bt_Iscan"
$ws1.Cells.Item(5, 2).Value = 0.4484
$ws1.Cells.Item(5, 3).Value = 11.0226
$ws1.Cells.Item(5, 4).Value = 0.0793
$ws1.Cells.Item(5, 5).Value = 1.8486
$ws1.Cells.Item(5, 6).Value = 100
$ws1.Cells.Item(6, 1).Value = "bt_PIscan"
$ws1.Cells.Item(6, 2).Value = 0.564
$ws1.Cells.Item(6, 3).Value = 14.8009
$ws1.Cells.Item(6, 4).Value = 0.0753
$ws1.Cells.Item(6, 5).Value = 2.5077
$ws1.Cells.Item(6, 6).Value = 100
$ws1.Cells.Item(7, 1).Value = "BT_ACL_Sniff_1dot28s_Master_0dBm"
$ws1.Cells.Item(7, 2).Value = 0.2675
$ws1.Cells.Item(7, 3).Value = 3.546
$ws1.Cells.Item(7, 4).Value = 0.0709
$ws1.Cells.Item(7, 5).Value = 0.6367
$ws1.Cells.Item(7, 6).Value = 100
$ws1.Cells.Item(8, 1).Value = "BT_ACL_Sniff_1dot28s_Master_4dBm"
$ws1.Cells.Item(8, 2).Value = 0.442
$ws1.Cells.Item(8, 3).Value = 4.5064
$ws1.Cells.Item(8, 4).Value = 0.0786
$ws1.Cells.Item(8, 5).Value = 1.0099
$ws1.Cells.Item(8, 6).Value = 100
$ws1.Cells.Item(9, 1).Value = "BT_ACL_Sniff_1dot28s_Master_12dot5dBm"
$ws1.Cells.Item(9, 2).Value = 0.3904
$ws1.Cells.Item(9, 3).Value = 5.9808
$ws1.Cells.Item(9, 4).Value = 0.0772
$ws1.Cells.Item(9, 5).Value = 1.1282
$ws1.Cells.Item(9, 6).Value = 100
$ws1.Cells.Item(10, 1).Value = "BT_ACL_Sniff_0dot5s_Master_0dBm"
$ws1.Cells.Item(10, 2).Value = 0.4586
$ws1.Cells.Item(10, 3).Value = 6.1246
$ws1.Cells.Item(10, 4).Value = 0.0636
$ws1.Cells.Item(10, 5).Value = 1.1492
$ws1.Cells.Item(10, 6).Value = 100
$ws1.Cells.Item(11, 1).Value = "BT_ACL_Sniff_0dot5s_Master_4dBm"
$ws1.Cells.Item(11, 2).Value = 0.3735
$ws1.Cells.Item(11, 3).Value = 6.0667
$ws1.Cells.Item(11, 4).Value = 0.0708
$ws1.Cells.Item(11, 5).Value = 1.0571
$ws1.Cells.Item(11, 6).Value = 100
$ws1.Cells.Item(12, 1).Value = "BT_SCO_HV3_Master_0dBm"
$ws1.Cells.Item(12, 2).Value = 8.1549
$ws1.Cells.Item(12, 3).Value = 9.4889
$ws1.Cells.Item(12, 4).Value = 7.8679
$ws1.Cells.Item(12, 5).Value = 0.2858
$ws1.Cells.Item(12, 6).Value = 100
$ws1.Cells.Item(13, 1).Value = "BT_SCO_HV3_Master_4dBm"
$ws1.Cells.Item(13, 2).Value = 8.5943
$ws1.Cells.Item(13, 3).Value = 9.8775
$ws1.Cells.Item(13, 4).Value = 8.1714
$ws1.Cells.Item(13, 5).Value = 0.3251
$ws1.Cells.Item(13, 6).Value = 100
$ws1.Cells.Item(14, 1).Value = "BT_SCO_HV3_Master_12dot5dBm"
$ws1.Cells.Item(14, 2).Value = 13.2779
$ws1.Cells.Item(14, 3).Value = 14.9539
$ws1.Cells.Item(14, 4).Value = 12.0784
$ws1.Cells.Item(14, 5).Value = 0.8209
$ws1.Cells.Item(14, 6).Value = 100
$ws1.Cells.Item(15, 1).Value = "BT_SCO_EV3_Master_0dBm"
$ws1.Cells.Item(15, 2).Value = 8.131
$ws1.Cells.Item(15, 3).Value = 9.3553
$ws1.Cells.Item(15, 4).Value = 7.8354
$ws1.Cells.Item(15, 5).Value = 0.2682
$ws1.Cells.Item(15, 6).Value = 100
$ws1.Cells.Item(16, 1).Value = "BT_SCO_EV3_Master_4dBm"
$ws1.Cells.Item(16, 2).Value = 8.6443
$ws1.Cells.Item(16, 3).Value = 9.4952
$ws1.Cells.Item(16, 4).Value = 8.1934
$ws1.Cells.Item(16, 5).Value = 0.3136
$ws1.Cells.Item(16, 6).Value = 100
$ws1.Cells.Item(17, 1).Value = "BT_SCO_EV3_Master_12dot5dBm"
$ws1.Cells.Item(17, 2).Value = 13.4639
$ws1.Cells.Item(17, 3).Value = 15.4356
$ws1.Cells.Item(17, 4).Value = 12.3822
$ws1.Cells.Item(17, 5).Value = 0.7764
$ws1.Cells.Item(17, 6).Value = 100
$ws1.Cells.Item(18, 1).Value = "BLE_Adv_1dot28s_3Channel_0dBm"
$ws1.Cells.Item(18, 2).Value = 0.1119
$ws1.Cells.Item(18, 3).Value = 1.2701
$ws1.Cells.Item(18, 4).Value = 0.0779
$ws1.Cells.Item(18, 5).Value = 0.1177
$ws1.Cells.Item(18, 6).Value = 100
$ws1.Cells.Item(19, 1).Value = "BLE_Adv_1dot28s_3Channel_4dBm"
$ws1.Cells.Item(19, 2).Value = 0.1006
$ws1.Cells.Item(19, 3).Value = 0.1796
$ws1.Cells.Item(19, 4).Value = 0.0765
$ws1.Cells.Item(19, 5).Value = 0.0183
$ws1.Cells.Item(19, 6).Value = 100
$ws1.Cells.Item(20, 1).Value = "BLE_Adv_1dot28s_3Channel_12dot5dBm"
$ws1.Cells.Item(20, 2).Value = 0.2259
$ws1.Cells.Item(20, 3).Value = 6.5985
$ws1.Cells.Item(20, 4).Value = 0.0765
$ws1.Cells.Item(20, 5).Value = 0.8808
$ws1.Cells.Item(20, 6).Value = 100
$ws1.Cells.Item(21, 1).Value = "BLE_Scan_1dot28s"
$ws1.Cells.Item(21, 2).Value = 0.3211
$ws1.Cells.Item(21, 3).Value = 11.2443
$ws1.Cells.Item(21, 4).Value = 0.0798
$ws1.Cells.Item(21, 5).Value = 1.5545
$ws1.Cells.Item(21, 6).Value = 100
$ws1.Cells.Item(22, 1).Value = "BLE_Scan_1s"
$ws1.Cells.Item(22, 2).Value = 0.2528
$ws1.Cells.Item(22, 3).Value = 11.2629
$ws1.Cells.Item(22, 4).Value = 0.0786
$ws1.Cells.Item(22, 5).Value = 1.1493
$ws1.Cells.Item(22, 6).Value = 100
$ws1.Cells.Item(23, 1).Value = "BLE_Scan_10ms"
$ws1.Cells.Item(23, 2).Value = 9.5751
$ws1.Cells.Item(23, 3).Value = 14.1172
$ws1.Cells.Item(23, 4).Value = 2.2144
$ws1.Cells.Item(23, 5).Value = 3.0326
$ws1.Cells.Item(23, 6).Value = 100
$ws1.Cells.Item(24, 1).Value = "BLE_Conn_1dot28s_0dBm"
$ws1.Cells.Item(24, 2).Value = 0.2144
$ws1.Cells.Item(24, 3).Value = 3.2109
$ws1.Cells.Item(24, 4).Value = 0.056
$ws1.Cells.Item(24, 5).Value = 0.5131
$ws1.Cells.Item(24, 6).Value = 100
$ws1.Cells.Item(25, 1).Value = "BLE_Conn_1dot28s_4dBm"
$ws1.Cells.Item(25, 2).Value = 0.1844
$ws1.Cells.Item(25, 3).Value = 2.7666
$ws1.Cells.Item(25, 4).Value = 0.0806
$ws1.Cells.Item(25, 5).Value = 0.4136
$ws1.Cells.Item(25, 6).Value = 100
$ws1.Cells.Item(26, 1).Value = "BLE_Conn_1dot28s_12dot5dBm"
$ws1.Cells.Item(26, 2).Value = 0.1867
$ws1.Cells.Item(26, 3).Value = 2.8814
$ws1.Cells.Item(26, 4).Value = 0.0608
$ws1.Cells.Item(26, 5).Value = 0.4286
$ws1.Cells.Item(26, 6).Value = 100

# The last row (idx 27) was a duplicate produced by the race; remove it.
$ws1.Rows.Item(27).Delete()

$ws2 = $wb.Worksheets.Item(2)
# --- Sheet "1_8": relabel rows and fix values ---
$ws2.Cells.Item(2, 1).Value = "Deep_Sleep"
$ws2.Cells.Item(2, 2).Value = 0.2344
$ws2.Cells.Item(2, 3).Value = 0.2444
$ws2.Cells.Item(2, 4).Value = 0.2211
$ws2.Cells.Item(2, 5).Value = 0.004
$ws2.Cells.Item(2, 6).Value = 100
$ws2.Cells.Item(3, 1).Value = "BT_Idle"
$ws2.Cells.Item(3, 2).Value = 8.328
$ws2.Cells.Item(3, 3).Value = 8.3428
$ws2.Cells.Item(3, 4).Value = 8.3072
$ws2.Cells.Item(3, 5).Value = 0.0096
$ws2.Cells.Item(3, 6).Value = 100
$ws2.Cells.Item(4, 1).Value = "BT_Pscan"
$ws2.Cells.Item(4, 2).Value = 0.3182
$ws2.Cells.Item(4, 3).Value = 3.2812
$ws2.Cells.Item(4, 4).Value = 0.2235
$ws2.Cells.Item(4, 5).Value = 0.4204
$ws2.Cells.Item(4, 6).Value = 100
$ws2.Cells.Item(5, 1).Value = "bt_Iscan"
$ws2.Cells.Item(5, 2).Value = 0.469
$ws2.Cells.Item(5, 3).Value = 6.1327
$ws2.Cells.Item(5, 4).Value = 0.22
$ws2.Cells.Item(5, 5).Value = 1.0239
$ws2.Cells.Item(5, 6).Value = 100
$ws2.Cells.Item(6, 1).Value = "bt_PIscan"
$ws2.Cells.Item(6, 2).Value = 0.5286
$ws2.Cells.Item(6, 3).Value = 8.2569
$ws2.Cells.Item(6, 4).Value = 0.2238
$ws2.Cells.Item(6, 5).Value = 1.4053
$ws2.Cells.Item(6, 6).Value = 100
$ws2.Cells.Item(7, 1).Value = "BT_ACL_Sniff_1dot28s_Master_0dBm"
$ws2.Cells.Item(7, 2).Value = 0.6602
$ws2.Cells.Item(7, 3).Value = 8.1711
$ws2.Cells.Item(7, 4).Value = 0.2224
$ws2.Cells.Item(7, 5).Value = 1.4835
$ws2.Cells.Item(7, 6).Value = 100
$ws2.Cells.Item(8, 1).Value = "BT_ACL_Sniff_1dot28s_Master_4dBm"
$ws2.Cells.Item(8, 2).Value = 0.562
$ws2.Cells.Item(8, 3).Value = 9.7773
$ws2.Cells.Item(8, 4).Value = 0.2233
$ws2.Cells.Item(8, 5).Value = 1.5503
$ws2.Cells.Item(8, 6).Value = 100
$ws2.Cells.Item(9, 1).Value = "BT_ACL_Sniff_1dot28s_Master_12dot5dBm"
$ws2.Cells.Item(9, 2).Value = 0.3384
$ws2.Cells.Item(9, 3).Value = 4.0969
$ws2.Cells.Item(9, 4).Value = 0.2244
$ws2.Cells.Item(9, 5).Value = 0.4439
$ws2.Cells.Item(9, 6).Value = 100
$ws2.Cells.Item(10, 1).Value = "BT_ACL_Sniff_0dot5s_Master_0dBm"
$ws2.Cells.Item(10, 2).Value = 0.5448
$ws2.Cells.Item(10, 3).Value = 7.8059
$ws2.Cells.Item(10, 4).Value = 0.2233
$ws2.Cells.Item(10, 5).Value = 1.2319
$ws2.Cells.Item(10, 6).Value = 100
$ws2.Cells.Item(11, 1).Value = "BT_ACL_Sniff_0dot5s_Master_4dBm"
$ws2.Cells.Item(11, 2).Value = 0.5593
$ws2.Cells.Item(11, 3).Value = 7.7425
$ws2.Cells.Item(11, 4).Value = 0.2208
$ws2.Cells.Item(11, 5).Value = 1.2682
$ws2.Cells.Item(11, 6).Value = 100
$ws2.Cells.Item(12, 1).Value = "BT_SCO_HV3_Master_0dBm"
$ws2.Cells.Item(12, 2).Value = 11.3549
$ws2.Cells.Item(12, 3).Value = 12.1485
$ws2.Cells.Item(12, 4).Value = 10.9978
$ws2.Cells.Item(12, 5).Value = 0.2769
$ws2.Cells.Item(12, 6).Value = 100
$ws2.Cells.Item(13, 1).Value = "BT_SCO_HV3_Master_4dBm"
$ws2.Cells.Item(13, 2).Value = 12.3031
$ws2.Cells.Item(13, 3).Value = 13.2164
$ws2.Cells.Item(13, 4).Value = 11.8003
$ws2.Cells.Item(13, 5).Value = 0.342
$ws2.Cells.Item(13, 6).Value = 100
$ws2.Cells.Item(14, 1).Value = "BT_SCO_HV3_Master_12dot5dBm"
$ws2.Cells.Item(14, 2).Value = 9.7492
$ws2.Cells.Item(14, 3).Value = 10.2441
$ws2.Cells.Item(14, 4).Value = 9.6097
$ws2.Cells.Item(14, 5).Value = 0.0886
$ws2.Cells.Item(14, 6).Value = 100
$ws2.Cells.Item(15, 1).Value = "BT_SCO_EV3_Master_0dBm"
$ws2.Cells.Item(15, 2).Value = 11.4917
$ws2.Cells.Item(15, 3).Value = 12.3767
$ws2.Cells.Item(15, 4).Value = 11.0183
$ws2.Cells.Item(15, 5).Value = 0.2732
$ws2.Cells.Item(15, 6).Value = 100
$ws2.Cells.Item(16, 1).Value = "BT_SCO_EV3_Master_4dBm"
$ws2.Cells.Item(16, 2).Value = 12.2369
$ws2.Cells.Item(16, 3).Value = 13.0403
$ws2.Cells.Item(16, 4).Value = 11.7453
$ws2.Cells.Item(16, 5).Value = 0.3438
$ws2.Cells.Item(16, 6).Value = 100
$ws2.Cells.Item(17, 1).Value = "BT_SCO_EV3_Master_12dot5dBm"
$ws2.Cells.Item(17, 2).Value = 9.7411
$ws2.Cells.Item(17, 3).Value = 10.2708
$ws2.Cells.Item(17, 4).Value = 9.5933
$ws2.Cells.Item(17, 5).Value = 0.1114
$ws2.Cells.Item(17, 6).Value = 100
$ws2.Cells.Item(18, 1).Value = "BLE_Adv_1dot28s_3Channel_0dBm"
$ws2.Cells.Item(18, 2).Value = 0.307
$ws2.Cells.Item(18, 3).Value = 3.9448
$ws2.Cells.Item(18, 4).Value = 0.2244
$ws2.Cells.Item(18, 5).Value = 0.5197
$ws2.Cells.Item(18, 6).Value = 100
$ws2.Cells.Item(19, 1).Value = "BLE_Adv_1dot28s_3Channel_4dBm"
$ws2.Cells.Item(19, 2).Value = 0.2754
$ws2.Cells.Item(19, 3).Value = 4.3856
$ws2.Cells.Item(19, 4).Value = 0.2172
$ws2.Cells.Item(19, 5).Value = 0.414
$ws2.Cells.Item(19, 6).Value = 100
$ws2.Cells.Item(20, 1).Value = "BLE_Adv_1dot28s_3Channel_12dot5dBm"
$ws2.Cells.Item(20, 2).Value = 0.2323
$ws2.Cells.Item(20, 3).Value = 0.3128
$ws2.Cells.Item(20, 4).Value = 0.218
$ws2.Cells.Item(20, 5).Value = 0.0089
$ws2.Cells.Item(20, 6).Value = 100
$ws2.Cells.Item(21, 1).Value = "BLE_Scan_1dot28s"
$ws2.Cells.Item(21, 2).Value = 0.2546
$ws2.Cells.Item(21, 3).Value = 1.7978
$ws2.Cells.Item(21, 4).Value = 0.2183
$ws2.Cells.Item(21, 5).Value = 0.176
$ws2.Cells.Item(21, 6).Value = 100
$ws2.Cells.Item(22, 1).Value = "BLE_Scan_1s"
$ws2.Cells.Item(22, 2).Value = 0.3905
$ws2.Cells.Item(22, 3).Value = 9.3187
$ws2.Cells.Item(22, 4).Value = 0.2214
$ws2.Cells.Item(22, 5).Value = 1.0629
$ws2.Cells.Item(22, 6).Value = 100
$ws2.Cells.Item(23, 1).Value = "BLE_Scan_10ms"
$ws2.Cells.Item(23, 2).Value = 8.0462
$ws2.Cells.Item(23, 3).Value = 12.7251
$ws2.Cells.Item(23, 4).Value = 0.8369
$ws2.Cells.Item(23, 5).Value = 2.8992
$ws2.Cells.Item(23, 6).Value = 100
$ws2.Cells.Item(24, 1).Value = "BLE_Conn_1dot28s_0dBm"
$ws2.Cells.Item(24, 2).Value = 0.2832
$ws2.Cells.Item(24, 3).Value = 2.7502
$ws2.Cells.Item(24, 4).Value = 0.2202
$ws2.Cells.Item(24, 5).Value = 0.3437
$ws2.Cells.Item(24, 6).Value = 100
$ws2.Cells.Item(25, 1).Value = "BLE_Conn_1dot28s_4dBm"
$ws2.Cells.Item(25, 2).Value = 0.2906
$ws2.Cells.Item(25, 3).Value = 3.2073
$ws2.Cells.Item(25, 4).Value = 0.2205
$ws2.Cells.Item(25, 5).Value = 0.4142
$ws2.Cells.Item(25, 6).Value = 100
$ws2.Cells.Item(26, 1).Value = "BLE_Conn_1dot28s_12dot5dBm"
$ws2.Cells.Item(26, 2).Value = 0.2806
$ws2.Cells.Item(26, 3).Value = 3.0164
$ws2.Cells.Item(26, 4).Value = 0.2169
$ws2.Cells.Item(26, 5).Value = 0.3526
$ws2.Cells.Item(26, 6).Value = 100

$ws2.Rows.Item(27).Delete()
